# "Pause switch renamed as continue switch"
#
# Slide 1 contains a single top-level group shape ("Group 19") that holds all
# the diagram's boxes, connectors and labels. The label for the GPIO 12
# pedal ("TextBox 13", shape id 14) currently reads "PAUSE pedal (GPIO 12)".
# It must become "CONTINUE pedal (GPIO 12)", and because the new word is
# wider than the old one the text box also needs to be widened to keep the
# autosized text box properly fitting its text (matching the target width).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$group = $s.Shapes.Item(1)

$targetShape = $null
for ($i = 1; $i -le $group.GroupItems.Count; $i++) {
    $item = $group.GroupItems.Item($i)
    if ($item.Id -eq 14) {
        $targetShape = $item
        break
    }
}

# Only change the first run ("PAUSE ") and leave the other runs
# ("pedal", " (GPIO 12)") untouched so their formatting is preserved.
$firstRun = $targetShape.TextFrame.TextRange.Characters(1, 6)
$firstRun.Text = "CONTINUE "

# Widen the textbox so the auto-fitted shape matches the new, longer caption.
$targetShape.Width = 169.6659055
